$d = $word.ActiveDocument

# --- Step 1: capture the "Meta description" paragraph (paragraph 2) before removing it ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Copy()

# --- Step 2: paste a copy of that paragraph right before the very last paragraph,
#     anchoring the insertion on a plain (unformatted) paragraph so nothing odd is inherited ---
$count = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($count - 1)
$anchor.Range.InsertParagraphAfter() | Out-Null

$newParaObj = $d.Paragraphs.Item($count)
$newParaObj.Style = "Normal"
$newParaObj.Range.Paste()

# --- Step 3: turn the pasted paragraph's text into the new bold heading-style line ---
$newRange = $d.Paragraphs.Item($count).Range
$newRange.Find.ClearFormatting()
$newRange.Find.Execute(
    "Meta description: Experience the Chinese-themed Double Happiness slot for free. Read our review on the game's features, symbols, and winning probability.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Play Double Happiness Slot for Free - Review 2021", 2) | Out-Null

# --- Step 4: update the (now second-to-last, originally last) "Prompt:" paragraph's text,
#     keeping its existing italic run formatting ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.Find.ClearFormatting()
$lastRange.Find.Execute(
    "Prompt: Please create a feature image for the slot game " + [char]34 + "Double Happiness" + [char]34 + " that fits the game's theme and features a happy Maya warrior with glasses. The image should be in cartoon style.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Experience the Chinese-themed Double Happiness slot for free. Read our review on the game's features, symbols, and winning probability.", 2) | Out-Null

# --- Step 5: remove the original "Meta description" paragraph near the top of the document ---
$d.Paragraphs.Item(2).Range.Delete()
